# Generate Report for Handback
#
# For the "ca3bfcdb-e0fa-4663-8b8e-c8a74059aee7" file, the handback report
# generator discovered that the handback file on hand is stale (there is a
# newer one). This fills in the previously-empty "Latest Target File",
# "Latest Handback File", "Latest Handback DateTime" and "Error Detail"
# columns for that row on both the zh-cn and de-de status sheets, widens the
# Error Detail column so the message is readable, and links the new target
# file name back to its source commit (mirroring the existing hyperlink on
# column A for the same row).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/120fdaff879f99672daf3adb96ac04f9e03b99cd/e2e/ca3bfcdb-e0fa-4663-8b8e-c8a74059aee7.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/ac8fbecd92f6f4fa4dce5b2a24e8cf1731d61509/e2e/ca3bfcdb-e0fa-4663-8b8e-c8a74059aee7.md."
$currentHandbackUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/120fdaff879f99672daf3adb96ac04f9e03b99cd/e2e/ca3bfcdb-e0fa-4663-8b8e-c8a74059aee7.md"
$displayName = "ca3bfcdb-e0fa-4663-8b8e-c8a74059aee7.md"

function Update-HandbackRow($ws, $targetXlf, $handbackDateTime) {
    # Error Detail column (P) is mostly empty; make it wide enough to read.
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # Latest Target File (I8): becomes a hyperlink to the current (stale)
    # handback commit, same visual style as the other file-name hyperlinks.
    $ws.Range("I8").Value = $displayName
    $ws.Range("I8").Font.Underline = $true
    $ws.Range("I8").Font.Color = 15570276
    $ws.Hyperlinks.Add($ws.Range("I8"), $currentHandbackUrl, [Type]::Missing, [Type]::Missing, $displayName) | Out-Null

    # Latest Handback File (J8)
    $ws.Range("J8").Value = $targetXlf

    # Latest Handback DateTime (K8) - stored as text like its siblings.
    $ws.Range("K8").Value = $handbackDateTime

    # Error Detail (P8)
    $ws.Range("P8").Value = $errorDetail
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow $wsZhCn "ca3bfcdb-e0fa-4663-8b8e-c8a74059aee7.783d65f37f017ca03f904efdeb51bbe78b795abc.zh-cn.xlf" "2016-08-13 18:53:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow $wsDeDe "ca3bfcdb-e0fa-4663-8b8e-c8a74059aee7.783d65f37f017ca03f904efdeb51bbe78b795abc.de-de.xlf" "2016-08-13 18:53:41"
